$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 17870744
$ws.Range("I62").Value = 22742174
$ws.Range("J62").Value = 8833.333000000001
$ws.Range("K62").Value = 22742174
$ws.Range("L62").Value = 8833.333000000001
$ws.Range("M62").Value = -22741550
$ws.Range("N62").Value = -10081.333

$ws.Range("H65").Value = 17870744
$ws.Range("I65").Value = 22742174
$ws.Range("J65").Value = 8833.333000000001
$ws.Range("K65").Value = 113710870
$ws.Range("L65").Value = 44166.665
$ws.Range("M65").Value = -113707750
$ws.Range("N65").Value = -50406.665

$ws.Range("H70").Value = 3619.25
$ws.Range("I70").Value = 3997
$ws.Range("J70").Value = 2486
$ws.Range("K70").Value = 11991
$ws.Range("L70").Value = 7458
$ws.Range("M70").Value = -11721
$ws.Range("N70").Value = -7998

$ws.Range("H73").Value = 3619.25
$ws.Range("I73").Value = 3997
$ws.Range("J73").Value = 2486
$ws.Range("K73").Value = 11991
$ws.Range("L73").Value = 7458
$ws.Range("M73").Value = -11055
$ws.Range("N73").Value = -9330

$ws.Range("H101").Value = 15485008
$ws.Range("I101").Value = 999999
$ws.Range("J101").Value = 22727512
$ws.Range("K101").Value = 2999997
$ws.Range("L101").Value = 68182536
$ws.Range("M101").Value = -2998375
$ws.Range("N101").Value = -68185780

$ws.Range("H132").Value = 41672916
$ws.Range("I132").Value = 34096284
$ws.Range("J132").Value = 125015880
$ws.Range("K132").Value = 102288852
$ws.Range("L132").Value = 375047640
$ws.Range("M132").Value = -102286322
$ws.Range("N132").Value = -375052700

$ws.Range("H138").Value = 3883.7856
$ws.Range("I138").Value = 2678.5264
$ws.Range("J138").Value = 4236.0923
$ws.Range("K138").Value = 8035.5792
$ws.Range("L138").Value = 12708.2769
$ws.Range("M138").Value = -2895.5792
$ws.Range("N138").Value = -22988.2769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14507800
$ws.Range("I32").Value = 20013540
$ws.Range("J32").Value = 19011.053
$ws.Range("K32").Value = 20013540
$ws.Range("L32").Value = 19011.053
$ws.Range("M32").Value = -20013253
$ws.Range("N32").Value = -19585.053

$ws.Range("H45").Value = 34874.734
$ws.Range("I45").Value = 41433.2
$ws.Range("J45").Value = 2082.4
$ws.Range("K45").Value = 41433.2
$ws.Range("L45").Value = 2082.4
$ws.Range("M45").Value = -41056.2
$ws.Range("N45").Value = -2836.4

$ws.Range("H122").Value = 4904.5884
$ws.Range("I122").Value = 5291.8667
$ws.Range("K122").Value = 15875.6001
$ws.Range("M122").Value = -13425.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 37560.625
$ws.Range("J109").Value = 37560.625
$ws.Range("L109").Value = 37560.625
$ws.Range("N109").Value = -40334.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2244.2273
$ws.Range("I31").Value = 1558.825
$ws.Range("K31").Value = 1558.825
$ws.Range("M31").Value = -1263.825

$ws.Range("H34").Value = 2244.2273
$ws.Range("I34").Value = 1558.825
$ws.Range("K34").Value = 1558.825
$ws.Range("M34").Value = -1356.825

$ws.Range("H41").Value = 8059
$ws.Range("I41").Value = 8059
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 8059
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -7631
$ws.Range("N41").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 12501143
$ws.Range("I22").Value = 50000750
$ws.Range("J22").Value = 1273.6666
$ws.Range("K22").Value = 150002250
$ws.Range("L22").Value = 3820.9998
$ws.Range("M22").Value = -150002081
$ws.Range("N22").Value = -4158.9998

$ws.Range("H26").Value = 575
$ws.Range("J26").Value = 2000
$ws.Range("L26").Value = 6000
$ws.Range("N26").Value = -6576

$ws.Range("H27").Value = 12501143
$ws.Range("I27").Value = 50000750
$ws.Range("J27").Value = 1273.6666
$ws.Range("K27").Value = 150002250
$ws.Range("L27").Value = 3820.9998
$ws.Range("M27").Value = -150002148
$ws.Range("N27").Value = -4024.9998

$ws.Range("H68").Value = 66674144
$ws.Range("I68").Value = 142857710
$ws.Range("J68").Value = 13521.125
$ws.Range("K68").Value = 428573130
$ws.Range("L68").Value = 40563.375
$ws.Range("M68").Value = -428572319
$ws.Range("N68").Value = -42185.375

$ws.Range("H71").Value = 66674144
$ws.Range("I71").Value = 142857710
$ws.Range("J71").Value = 13521.125
$ws.Range("K71").Value = 1285719390
$ws.Range("L71").Value = 121690.125
$ws.Range("M71").Value = -1285715334
$ws.Range("N71").Value = -129802.125

$ws.Range("H121").Value = 5129081
$ws.Range("I121").Value = 386
$ws.Range("J121").Value = 8334515
$ws.Range("K121").Value = 1158
$ws.Range("L121").Value = 25003545
$ws.Range("M121").Value = 152
$ws.Range("N121").Value = -25006165

$ws.Range("H122").Value = 31257170
$ws.Range("I122").Value = 166667090
$ws.Range("J122").Value = 8726.846
$ws.Range("K122").Value = 1500003810
$ws.Range("L122").Value = 78541.614
$ws.Range("M122").Value = -1500001360
$ws.Range("N122").Value = -83441.614

$ws.Range("H131").Value = 735.91
$ws.Range("I131").Value = 470
$ws.Range("J131").Value = 759.0326
$ws.Range("K131").Value = 1410
$ws.Range("L131").Value = 2277.0978
$ws.Range("M131").Value = 3630
$ws.Range("N131").Value = -12357.0978

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5250047.5
$ws.Range("I12").Value = 5250047.5
$ws.Range("K12").Value = 5250047.5
$ws.Range("M12").Value = -5249907.5

$ws.Range("H80").Value = 5266521
$ws.Range("I80").Value = 3992
$ws.Range("J80").Value = 14287999
$ws.Range("K80").Value = 3992
$ws.Range("L80").Value = 14287999
$ws.Range("M80").Value = -2994
$ws.Range("N80").Value = -14289995

$ws.Range("H83").Value = 5266521
$ws.Range("I83").Value = 3992
$ws.Range("J83").Value = 14287999
$ws.Range("K83").Value = 19960
$ws.Range("L83").Value = 71439995
$ws.Range("M83").Value = -14968
$ws.Range("N83").Value = -71449979

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1613.909
$ws.Range("I68").Value = 1594.125
$ws.Range("J68").Value = 1666.6666
$ws.Range("K68").Value = 1594.125
$ws.Range("L68").Value = 1666.6666
$ws.Range("M68").Value = -845.125
$ws.Range("N68").Value = -3164.6666

$ws.Range("H71").Value = 1613.909
$ws.Range("I71").Value = 1594.125
$ws.Range("J71").Value = 1666.6666
$ws.Range("K71").Value = 7970.625
$ws.Range("L71").Value = 8333.333000000001
$ws.Range("M71").Value = -4226.625
$ws.Range("N71").Value = -15821.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4376
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -21880
$ws.Range("N65").ClearContents()

$ws.Range("H74").Value = 11400
$ws.Range("J74").Value = 9800
$ws.Range("L74").Value = 9800
$ws.Range("N74").Value = -11672

$ws.Range("H77").Value = 11400
$ws.Range("J77").Value = 9800
$ws.Range("L77").Value = 29400
$ws.Range("N77").Value = -38760
